# Adding duplicates removal script
# Append two more bibliography entries to the "sources" table (Tabla1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Δέφνερ, Μιχαήλ (1923) — no URL
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Δέφνερ, Μιχαήλ (1923). Λεξικόν της Τσακώνικης Διαλέκτου"

# Row 4: Warr, John — Tsakoniandialect.info, with URL
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Warr, John. Tsakoniandialect.info"
$ws.Range("C4").Value = "http://www.tsakoniandialect.info/"

# Grow the table ("Tabla1") so the new rows are recognised as part of it
$tbl = $ws.ListObjects("Tabla1")
$tbl.Resize($ws.Range("A1:C4"))

# Match the active selection left behind in the edited workbook
$ws.Range("B3").Select()

$wb.Save()
